$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.645.26"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").Value = "2.384.66"
$ws.Range("E3").Value = "  +3.29%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'309.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").Value = "'104.53"
$ws.Range("D6").ClearFormats()

$ws.Range("D7").Value = "'0.509"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.30%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.521"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.09%  "

$ws.Range("D10").Value = "'36.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").Value = "'53.41"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.09%  "

$ws.Range("D12").Value = "'0.0814"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").Value = "'7.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.22%  "

$ws.Range("D15").Value = "2.751.61"
$ws.Range("E15").Value = "  +3.22%  "

$ws.Range("D16").Value = "'15.66"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.19%  "

$ws.Range("D17").Value = "2.377.70"
$ws.Range("E17").Value = "  +3.31%  "

$ws.Range("D18").Value = "'0.812"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").Value = "43.606.58"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("D20").Value = "'6.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.74%  "

$ws.Range("D21").Value = "'11.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.14%  "

$ws.Range("D22").Value = "0.0₃0918"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").Value = "'68.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").Value = "'241.37"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("E25").Value = "  +1.96%  "

$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").Value = "'25.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.19%  "

$ws.Range("E29").Value = "  -3.25%  "

$ws.Range("D30").Value = "'36.59"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.55%  "

$ws.Range("D31").Value = "'9.54"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("D32").Value = "'2.10"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").Value = "'160.84"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.84%  "

$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").Value = "'18.34"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.43%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'3.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'2.53"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.60%  "

$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").Value = "'4.66"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.29%  "

$ws.Range("E41").Value = "  +5.74%  "

$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("D44").Value = "'2.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +12.82%  "

$ws.Range("D45").Value = "2.035.23"
$ws.Range("E45").Value = "  +2.40%  "

$ws.Range("D46").Value = "'19.63"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.87%  "

$ws.Range("D47").Value = "'0.0290"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").Value = "'3.14"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.55%  "

$ws.Range("D49").Value = "'10.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.42%  "

$ws.Range("D50").Value = "'57.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.89%  "

$ws.Range("D51").Value = "'2.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.24%  "
